# Rebuild Sheet1 as a 3-box UML class diagram (Room / Songs / Guests,
# plus corresponding test-double boxes) with Class:/Attribute:/Method:
# labels in column A.

function Set-GrayFill($range) {
    $range.Interior.Color = 14277081
}

function Set-Empty($range) {
    # Touch (but don't visibly format) a cell so it is materialised in the
    # sheet as an explicit-but-empty <c> element, matching cells that carry
    # a purely cosmetic "no border applied" flag in the source file.
    $range.Borders.LineStyle = -4142
}

function Set-Box($range, $left, $right, $top, $bottom) {
    Set-GrayFill $range
    # Start from all four thin sides, then strip the ones that shouldn't be
    # there - this avoids leaving unused intermediate border/style entries
    # behind (each distinct combination touched gets its own style slot).
    $range.Borders.LineStyle = 1
    $range.Borders.Weight = 2
    if (-not $left)   { $range.Borders.Item(7).LineStyle = -4142 }
    if (-not $right)  { $range.Borders.Item(10).LineStyle = -4142 }
    if (-not $top)    { $range.Borders.Item(8).LineStyle = -4142 }
    if (-not $bottom) { $range.Borders.Item(9).LineStyle = -4142 }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - the whole sheet content/layout is being replaced.
$ws.Cells.Clear()

$ws.Columns("A").ColumnWidth = 10

# ---- Column A labels ----------------------------------------------------
$ws.Range("A2").Value = "Class:"
$ws.Range("A4").Value = "Attribute:"
$ws.Range("A9").Value = "Method:"

# ---- Row 2: class name headers (full box) -------------------------------
$ws.Range("B2").Value = "room"
$ws.Range("D2").Value = "songs"
$ws.Range("F2").Value = "guests"
Set-Box $ws.Range("B2") $true $true $true $true
Set-Box $ws.Range("D2") $true $true $true $true
Set-Box $ws.Range("F2") $true $true $true $true

# ---- Row 3: blank spacer (filled, no border) -----------------------------
Set-GrayFill $ws.Range("B3")
Set-GrayFill $ws.Range("D3")
Set-GrayFill $ws.Range("F3")

# ---- Rows 4-7: attribute boxes -------------------------------------------
$ws.Range("B4").Value = "name"
$ws.Range("D4").Value = "name"
$ws.Range("F4").Value = "name"
Set-Box $ws.Range("B4") $true $true $true $false
Set-Box $ws.Range("D4") $true $true $true $false
Set-Box $ws.Range("F4") $true $true $true $false

$ws.Range("B5").Value = "type"
$ws.Range("D5").Value = "artist"
$ws.Range("F5").Value = "age"
Set-Box $ws.Range("B5") $true $true $false $false
Set-Box $ws.Range("D5") $true $true $false $false
Set-Box $ws.Range("F5") $true $true $false $false

$ws.Range("B6").Value = "capacity"
$ws.Range("D6").Value = "length"
$ws.Range("F6").Value = "money"
Set-Box $ws.Range("B6") $true $true $false $false
Set-Box $ws.Range("D6") $true $true $false $false
Set-Box $ws.Range("F6") $true $true $false $false

$ws.Range("D7").Value = "genre"
Set-Box $ws.Range("B7") $true $true $false $true
Set-Box $ws.Range("D7") $true $true $false $true
Set-Box $ws.Range("F7") $true $true $false $true

# Spacer columns C / E alongside the attribute box (rows 4-7) stay
# unfilled / borderless, but are still present as explicit empty cells.
"C4","E4","C5","E5","C6","E6","C7","E7" | ForEach-Object { Set-Empty $ws.Range($_) }

# ---- Row 8: blank spacer (filled, no border) ------------------------------
Set-GrayFill $ws.Range("B8")
Set-GrayFill $ws.Range("D8")
Set-GrayFill $ws.Range("F8")

# ---- Rows 9-12: method boxes ----------------------------------------------
$ws.Range("B9").Value = "test name"
$ws.Range("D9").Value = "test name"
$ws.Range("F9").Value = "test name"
Set-Box $ws.Range("B9") $true $true $true $false
Set-Box $ws.Range("D9") $true $true $true $false
Set-Box $ws.Range("F9") $true $true $true $false

$ws.Range("B10").Value = "test type"
$ws.Range("D10").Value = "test artist"
$ws.Range("F10").Value = "test age"
Set-Box $ws.Range("B10") $true $true $false $false
Set-Box $ws.Range("D10") $true $true $false $false
Set-Box $ws.Range("F10") $true $true $false $false

$ws.Range("B11").Value = "test capacity"
$ws.Range("D11").Value = "test length"
$ws.Range("F11").Value = "test money"
Set-Box $ws.Range("B11") $true $true $false $false
Set-Box $ws.Range("D11") $true $true $false $false
Set-Box $ws.Range("F11") $true $true $false $true

$ws.Range("B12").Value = "add_songs"
$ws.Range("D12").Value = "test genre"
Set-Box $ws.Range("B12") $true $true $false $false
Set-Box $ws.Range("D12") $true $true $false $true

$ws.Range("B13").Value = "check_in_guests"
Set-Box $ws.Range("B13") $true $true $false $false

$ws.Range("B14").Value = "check_out_guests"
Set-Box $ws.Range("B14") $true $true $false $true

# Spacer columns C / E alongside the method box (rows 9-12), plus F12 which
# also has no text (bottom-right corner of the guests box is shorter), stay
# unfilled / borderless, but are still present as explicit empty cells.
"C9","E9","C10","E10","C11","E11","C12","E12","F12" | ForEach-Object { Set-Empty $ws.Range($_) }

$ws.Range("B19").Select()
